# Applies the 'Updated cryptos list' data refresh: new Price (column D)
# and Volume(1h) (column E) values for rows 2-51 (row 41 unchanged).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D sometimes holds purely-numeric-looking text (e.g. "304.49").
# Force those cells to keep their original Text format so COM does not
# auto-convert the assigned string into a numeric value.
function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
}

$ws.Range("D2").Value = "42.112.91"
$ws.Range("E2").Value = "  +0.62%  "
$ws.Range("D3").Value = "2.279.66"
$ws.Range("E3").Value = "  +0.26%  "
$ws.Range("E4").Value = "  +0.05%  "
Set-TextValue $ws.Range("D5") "154.84"
$ws.Range("E5").Value = "  +15,377.92%  "
Set-TextValue $ws.Range("D6") "304.49"
$ws.Range("E6").Value = "  +0.21%  "
Set-TextValue $ws.Range("D7") "93.93"
$ws.Range("E7").Value = "  +0.58%  "
$ws.Range("E8").Value = "  -0.16%  "
$ws.Range("E9").Value = "  -0.02%  "
Set-TextValue $ws.Range("D10") "0.490"
$ws.Range("E10").Value = "  +0.65%  "
Set-TextValue $ws.Range("D11") "34.00"
$ws.Range("E11").Value = "  +4.04%  "
Set-TextValue $ws.Range("D12") "0.0804"
$ws.Range("E12").Value = "  +0.24%  "
$ws.Range("E13").Value = "  -2.40%  "
Set-TextValue $ws.Range("D14") "6.67"
$ws.Range("E14").Value = "  -0.37%  "
$ws.Range("D15").Value = "2.634.11"
$ws.Range("E15").Value = "  +0.33%  "
Set-TextValue $ws.Range("D16") "14.34"
$ws.Range("E16").Value = "  +0.45%  "
$ws.Range("D17").Value = "2.277.45"
$ws.Range("E17").Value = "  +0.00%  "
$ws.Range("E18").Value = "  +3.68%  "
$ws.Range("D19").Value = "42.035.21"
$ws.Range("E19").Value = "  +0.58%  "
Set-TextValue $ws.Range("D20") "12.79"
$ws.Range("E20").Value = "  +4.20%  "
$ws.Range("D21").Value = "0.0₃0917"
$ws.Range("E21").Value = "  +1.00%  "
Set-TextValue $ws.Range("D22") "5.99"
$ws.Range("E22").Value = "  +0.67%  "
Set-TextValue $ws.Range("D23") "68.05"
$ws.Range("E23").Value = "  +1.05%  "
Set-TextValue $ws.Range("D24") "243.65"
$ws.Range("E24").Value = "  +0.07%  "
Set-TextValue $ws.Range("D25") "2.59"
$ws.Range("E25").Value = "  +0.69%  "
$ws.Range("E26").Value = "  +1.49%  "
$ws.Range("E27").Value = "  -0.06%  "
Set-TextValue $ws.Range("D28") "24.07"
$ws.Range("E28").Value = "  -1.05%  "
$ws.Range("E29").Value = "  +5.73%  "
Set-TextValue $ws.Range("D30") "9.67"
$ws.Range("E30").Value = "  +0.46%  "
$ws.Range("E31").Value = "  +1.15%  "
Set-TextValue $ws.Range("D32") "160.39"
$ws.Range("E32").Value = "  +1.24%  "
Set-TextValue $ws.Range("D33") "5.35"
$ws.Range("E33").Value = "  +2.95%  "
$ws.Range("E34").Value = "  +0.08%  "
Set-TextValue $ws.Range("D35") "0.0753"
$ws.Range("E35").Value = "  +0.21%  "
Set-TextValue $ws.Range("D36") "3.08"
$ws.Range("E36").Value = "  +0.74%  "
$ws.Range("E37").Value = "  +2.98%  "
Set-TextValue $ws.Range("D38") "17.01"
$ws.Range("E38").Value = "  +1.54%  "
$ws.Range("E39").Value = "  -0.39%  "
$ws.Range("E40").Value = "  -0.36%  "
Set-TextValue $ws.Range("D42") "4.20"
$ws.Range("E42").Value = "  +6.72%  "
Set-TextValue $ws.Range("D43") "19.81"
$ws.Range("E43").Value = "  +1.52%  "
$ws.Range("D44").Value = "2.022.30"
$ws.Range("E44").Value = "  -3.02%  "
$ws.Range("E45").Value = "  +11.32%  "
$ws.Range("E46").Value = "  +1.27%  "
Set-TextValue $ws.Range("D47") "10.22"
$ws.Range("E47").Value = "  -1.19%  "
Set-TextValue $ws.Range("D48") "2.93"
$ws.Range("E48").Value = "  +0.36%  "
Set-TextValue $ws.Range("D49") "53.53"
$ws.Range("E49").Value = "  +3.13%  "
$ws.Range("E50").Value = "  -1.18%  "
Set-TextValue $ws.Range("D51") "72.24"
$ws.Range("E51").Value = "  -1.14%  "
